$p = $ppt.ActivePresentation

# ------------------------------------------------------------------
# 1) Bump the auto "datetimeFigureOut" date placeholders on the slide
#    master and every custom (slide) layout from 3/30/2016 to
#    3/31/2016 (the deck was re-saved a day later).
# ------------------------------------------------------------------
function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shape = $shapes.Item($i)
        $phType = -1
        try { $phType = $shape.PlaceholderFormat.Type } catch { $phType = -1 }
        if ($phType -eq 16) {
            if ($shape.TextFrame.TextRange.Text -eq "3/30/2016") {
                $shape.TextFrame.TextRange.Text = "3/31/2016"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

# ------------------------------------------------------------------
# 2) Fix the estimated-hours figure on the poster's "Plan" panel:
#    "Total Estimated Hours Until Completion: 23" -> "... : 416"
# ------------------------------------------------------------------
$slide = $p.Slides.Item(1)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shape = $slide.Shapes.Item($i)
    if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
        $tr = $shape.TextFrame.TextRange
        $full = $tr.Text
        if ($full.IndexOf("Total Estimated Hours Until Completion: 23") -ge 0) {
            $idx0 = $full.IndexOf("Total Estimated Hours Until Completion: 23")

            $startCompletion = $idx0 + 1
            $lenCompletion = ("Total Estimated Hours Until Completion").Length
            $startColon = $startCompletion + $lenCompletion
            $lenColon = 2
            $startNum = $startColon + $lenColon
            $lenNum = 2

            $numRange = $tr.Characters($startNum, $lenNum)
            $numRange.Text = "416"

            $colonRange = $tr.Characters($startColon, $lenColon)
            $colonRange.Text = ": "

            $labelRange = $tr.Characters($startCompletion, $lenCompletion)
            $labelRange.Text = "Total Estimated Hours Until Completion"
        }
    }
}
